$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CustomerCreated-Event")
$ws.Activate()

# Rename the "ProtoBuffMessageType" message-type label used in row 3 to
# "ProtobufType" for the 1.6.3 release.
$ws.Range("K3").Value = "ProtobufType"

# Leave the selection where the edit was made.
$ws.Range("K3").Select() | Out-Null
